# 20 valid Paranthesis solution
# Adds a new row (row 8) to Sheet1 documenting LeetCode #20 "Valid Parentheses",
# marks it (and the two rows above it) with a "Diff" of "Easy", and widens the
# Category/Patterns columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new row-8 text cells first, and in the same order the author typed
# them (Suggestion, then Question name, then Category) so that the workbook's
# shared-string table grows in that order.
$ws.Range("H8").Value = "22. Generate Parentheses, 32. Longest Valid Parentheses, 301. Remove Invalid Parentheses"
$ws.Range("H8").WrapText = $true

$ws.Range("D8").Value = "Valid Paranthesis"
$ws.Range("E8").Value = "logic and concept and clarity"

# Remaining numeric / lookup cells for the new row.
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 20
$ws.Range("G8").Value = "Easy"

# Backfill the "Diff" column for the two existing rows just above (they were
# missing a difficulty rating).
$ws.Range("G6").Value = "Easy"
$ws.Range("G7").Value = "Easy"

# Re-fit the Category (E) and Patterns (F) columns now that they hold longer
# text.
$ws.Columns.Item(5).ColumnWidth = 23.6
$ws.Columns.Item(6).ColumnWidth = 16

# Leave the selection where the author ended up after entering the new row.
$ws.Range("H11").Select() | Out-Null
